$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "System" label in column A was originally repeated on every data row
# within each group (sharing the same style + shared string). Excel only
# needs the label on the first row of each group, so clear the redundant
# label cells from every row after the first one in each group, leaving the
# cell style intact but removing the text value.
$labelRange = $ws.Range("A3:A37,A40:A54,A57:A76,A79:A119,A122:A154,A157:A166,A169:A183,A186:A207")
foreach ($area in $labelRange.Areas) {
    $area.ClearContents()
}

# Restore the view to the top of the sheet with L2 selected (matches the
# saved workbook view state after the edit).
$ws.Range("L2").Select()
